$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Mac OS X" hint to the existing Linux references (3 distinct strings,
# used across 5 cells total thanks to shared-string de-duplication).
$null = $ws.Cells.Replace("start.sh - Linux)", "start.sh - Linux, Mac OS X)")
$null = $ws.Cells.Replace("similar - Linux)", "similar - Linux, Mac OS X)")
$null = $ws.Cells.Replace("top on Linux)", "top on Linux, Mac OS X)")

# New objective hint cell next to the existing "Objective" text (row 2).
$c2 = $ws.Range("C2")
$c2.Value = "Execute tests on windows (7, 8), Linux (e.g. Ubuntu) and Mac OS X"
$c2.WrapText = $true
$c2.Font.Bold = $true
$c2.Font.Size = 18

# Row 2 needs to grow to fit the new, taller content.
$ws.Rows.Item(2).RowHeight = 47.25

# Update the view: clear the old scrolled-down position and select C8.
$ws.Range("C8").Select()
